$d = $word.ActiveDocument

# Remove the Russian language formatting from paragraph 2 (both paragraph mark rPr and the "!" run rPr)
$para2 = $d.Paragraphs(2)
$para2.Range.Font.Reset()

# Add two new paragraphs at the end of the document (after "Are right here!!")
$endRange = $d.Paragraphs(2).Range
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertAfter("Third and forth")

$d.Paragraphs(3).Range.InsertParagraphAfter()
$d.Paragraphs(4).Range.InsertAfter("Lines are here!")
